$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.709.31"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.005.33"
$ws.Range("E3").Value = "  +3.27%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.42"
$ws.Range("E5").Value = "  +5.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.56"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.99"
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0847"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.90"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "3.477.27"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "3.007.55"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.978"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "51.703.29"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.45"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.76"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("E27").Value = "  +19.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.54"
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.25"
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.01"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.07"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.16"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0449"
$ws.Range("E36").Value = "  +5.90%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.64"
$ws.Range("E39").Value = "  +3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  -5.52%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.62"
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.43"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.08"
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("E46").Value = "  +7.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.276"
$ws.Range("E47").Value = "  +16.20%  "
$ws.Range("D48").Value = "2.060.68"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.29"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +13.87%  "
$ws.Range("E51").Value = "  +3.23%  "
